# Removed fetal deaths before PNC encounter.
$wb = $excel.ActiveWorkbook

$wsUntrt = $wb.Worksheets.Item("potential_preg_untrt")

# Set the fetal-death probabilities (column C, rows 2-8) to 0.
# Downstream formulas (E2:E8 on this sheet, C43:E49 / C84:E90 on this
# sheet, and the corresponding cells on potential_preg_trt) recalc
# automatically because they reference these cells.
$wsUntrt.Range("C2:C8").Value = 0

# Update which sheet/cell is active, matching the saved view state.
$wsUntrt.Range("C9").Select()
$wsUntrt.Activate()
